{"js": "const body = context.document.body;\n\n// 1) Title line: the meeting was actually held on a Monday, not a\n//    Thursday - fix the weekday name. The day number (\"26\") is already\n//    correct and stays untouched. The title is the very first paragraph;\n//    scope the search to it so the unrelated \"Donnerstag\" further down\n//    the document (in \"N\u00e4chste Zwischenbesprechung? Donnerstag,\n//    29.03.2018...\") is left untouched.\nconst titleParagraph = body.paragraphs.getFirst();\nconst weekdayMatch = titleParagraph.search(\"Donnerstag\", { matchCase: true });\nweekdayMatch.load(\"text\");\nawait context.sync();\nif (weekdayMatch.items.length > 0) {\n  weekdayMatch.items[0].insertText(\"Montag\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2) Word maintains a single \"_GoBack\" bookmark marking the location of\n//    the most recent edit. It currently sits right after the title's\n//    \"26\" - remove it here; it will be re-created after the edit below,\n//    which is the new \"most recent\" edit location.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 3) Closing note: this protocol was actually finished/sent on the 26th,\n//    not the 23rd - fix the typo in the date. Scope the replace to just\n//    the \"23\" digits so the \"Montag, \" / \".03\" text around it is left\n//    alone.\nconst noteMatch = body.search(\"Montag, 23.03\", { matchCase: true });\nnoteMatch.load(\"text\");\nawait context.sync();\nif (noteMatch.items.length > 0) {\n  const noteRange = noteMatch.items[0];\n  const dayDigits = noteRange.search(\"23\", { matchCase: true });\n  dayDigits.load(\"text\");\n  await context.sync();\n  if (dayDigits.items.length > 0) {\n    dayDigits.items[0].insertText(\"26\", Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n\n// 4) Re-create the \"_GoBack\" bookmark right after the corrected day\n//    number, mirroring Word's behaviour of tracking the latest edit.\nconst fixedNoteMatch = body.search(\"Montag, 26.03\", { matchCase: true });\nfixedNoteMatch.load(\"text\");\nawait context.sync();\nif (fixedNoteMatch.items.length > 0) {\n  const fixedNoteRange = fixedNoteMatch.items[0];\n  const fixedDayDigits = fixedNoteRange.search(\"26\", { matchCase: true });\n  fixedDayDigits.load(\"text\");\n  await context.sync();\n  if (fixedDayDigits.items.length > 0) {\n    const afterDay = fixedDayDigits.items[0].getRange(\"End\");\n    afterDay.insertBookmark(\"_GoBack\");\n    await context.sync();\n  }\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Title line: the meeting was actually held on a Monday, not a\n#    Thursday - fix the weekday name. The day number (\"26\") is already\n#    correct and stays untouched. Scope the Find to the title paragraph\n#    (the very first paragraph) so the unrelated \"Donnerstag\" further\n#    down the document (\"N\u00e4chste Zwischenbesprechung? Donnerstag,\n#    29.03.2018...\") is left untouched.\n$titleRange = $d.Paragraphs(1).Range\n$titleFind = $titleRange.Find\n$titleFind.Text = \"Donnerstag\"\n$titleFind.Replacement.Text = \"Montag\"\n$titleFind.Wrap = 0\n$titleFind.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2)\n\n# 2) Word maintains a single \"_GoBack\" bookmark marking the location of\n#    the most recent edit. It currently sits right after the title's\n#    \"26\" - remove it here; it gets re-created below at the new \"most\n#    recent\" edit location.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# 3) Closing note: this protocol was actually finished/sent on the 26th,\n#    not the 23rd - fix the typo in the date. Scope the Find to the\n#    closing-note paragraph and to just the \"23\" digits so the\n#    surrounding \"Montag, \" / \".03.2018\" text is left alone.\n$noteParagraph = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs($i).Range.Text -like \"*Dieses Dokument wurde am*\") {\n        $noteParagraph = $i\n        break\n    }\n}\n\n$noteRange = $d.Paragraphs($noteParagraph).Range\n$noteFind = $noteRange.Find\n$noteFind.Text = \"23\"\n$noteFind.Wrap = 0\n$noteFind.Execute()\n$noteRange.Text = \"26\"\n\n# 4) Re-create the \"_GoBack\" bookmark right after the corrected day\n#    number, mirroring Word's behaviour of tracking the latest edit.\n$noteRange.Collapse(0)\n$d.Bookmarks.Add(\"_GoBack\", $noteRange)\n"}
